$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation for "Ajo" (Chino / Primera) was recorded.
# It belongs chronologically at row 40, so insert a new row there and push
# every existing row from 40 downward (old row 40 -> new row 41, ...,
# old row 123 -> new row 124).
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new observation.
$ws.Range("A40").Value = 11
$ws.Range("B40").Value = "Vega Monumental Concepción"
$ws.Range("C40").Value = "Bíobío"
$ws.Range("D40").Value = 44530
$ws.Range("E40").Value = 8
$ws.Range("F40").Value = 100112003
$ws.Range("G40").Value = "Ajo"
$ws.Range("H40").Value = "Chino"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 310
$ws.Range("K40").Value = 21000
$ws.Range("L40").Value = 22000
$ws.Range("M40").Value = 21484
$ws.Range("N40").Value = "`$/caja 10 kilos"
$ws.Range("O40").Value = "China"
$ws.Range("P40").Value = 2148
$ws.Range("Q40").Value = 10
$ws.Range("R40").Value = "Hortaliza"
